$d = $word.ActiveDocument

$d.Content.Find.Execute("99-12=87", $true, $false, $false, $false, $false, $true, 1, $false, "80-30=50", 2) | Out-Null
$d.Content.Find.Execute("50-17=33", $true, $false, $false, $false, $false, $true, 1, $false, "71+24=95", 2) | Out-Null
$d.Content.Find.Execute("48-32=16", $true, $false, $false, $false, $false, $true, 1, $false, "45-5=40", 2) | Out-Null
$d.Content.Find.Execute("91-48=43", $true, $false, $false, $false, $false, $true, 1, $false, "43-36=7", 2) | Out-Null
$d.Content.Find.Execute("83-75=8", $true, $false, $false, $false, $false, $true, 1, $false, "12+45=57", 2) | Out-Null
$d.Content.Find.Execute("40-1=39", $true, $false, $false, $false, $false, $true, 1, $false, "11+39=50", 2) | Out-Null
$d.Content.Find.Execute("35+25=60", $true, $false, $false, $false, $false, $true, 1, $false, "43-24=19", 2) | Out-Null
$d.Content.Find.Execute("41-8=33", $true, $false, $false, $false, $false, $true, 1, $false, "98-76=22", 2) | Out-Null
$d.Content.Find.Execute("29-18=11", $true, $false, $false, $false, $false, $true, 1, $false, "86-69=17", 2) | Out-Null
$d.Content.Find.Execute("68+31=99", $true, $false, $false, $false, $false, $true, 1, $false, "21+64=85", 2) | Out-Null
$d.Content.Find.Execute("93-62=31", $true, $false, $false, $false, $false, $true, 1, $false, "20-8=12", 2) | Out-Null
$d.Content.Find.Execute("73+8=81", $true, $false, $false, $false, $false, $true, 1, $false, "68-57=11", 2) | Out-Null
$d.Content.Find.Execute("94-51=43", $true, $false, $false, $false, $false, $true, 1, $false, "90-14=76", 2) | Out-Null
$d.Content.Find.Execute("97-5=92", $true, $false, $false, $false, $false, $true, 1, $false, "92-4=88", 2) | Out-Null
$d.Content.Find.Execute("17-6=11", $true, $false, $false, $false, $false, $true, 1, $false, "0+34=34", 2) | Out-Null
$d.Content.Find.Execute("38+32=70", $true, $false, $false, $false, $false, $true, 1, $false, "12+64=76", 2) | Out-Null
$d.Content.Find.Execute("72+5=77", $true, $false, $false, $false, $false, $true, 1, $false, "2+79=81", 2) | Out-Null
$d.Content.Find.Execute("42-41=1", $true, $false, $false, $false, $false, $true, 1, $false, "8+55=63", 2) | Out-Null
$d.Content.Find.Execute("66-49=17", $true, $false, $false, $false, $false, $true, 1, $false, "4+54=58", 2) | Out-Null
$d.Content.Find.Execute("50+17=67", $true, $false, $false, $false, $false, $true, 1, $false, "59+30=89", 2) | Out-Null
$d.Content.Find.Execute("76-25=51", $true, $false, $false, $false, $false, $true, 1, $false, "19+29=48", 2) | Out-Null
$d.Content.Find.Execute("90-12=78", $true, $false, $false, $false, $false, $true, 1, $false, "6+34=40", 2) | Out-Null
$d.Content.Find.Execute("85+0=85", $true, $false, $false, $false, $false, $true, 1, $false, "63-43=20", 2) | Out-Null
$d.Content.Find.Execute("2+10=12", $true, $false, $false, $false, $false, $true, 1, $false, "67-24=43", 2) | Out-Null
$d.Content.Find.Execute("31+41=72", $true, $false, $false, $false, $false, $true, 1, $false, "23+1=24", 2) | Out-Null
$d.Content.Find.Execute("22+37=59", $true, $false, $false, $false, $false, $true, 1, $false, "73-66=7", 2) | Out-Null
$d.Content.Find.Execute("14+69=83", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("80-14=66", $true, $false, $false, $false, $false, $true, 1, $false, "11+77=88", 2) | Out-Null
$d.Content.Find.Execute("94-47=47", $true, $false, $false, $false, $false, $true, 1, $false, "48-23=25", 2) | Out-Null
$d.Content.Find.Execute("67+25=92", $true, $false, $false, $false, $false, $true, 1, $false, "46+38=84", 2) | Out-Null
$d.Content.Find.Execute("63-49=14", $true, $false, $false, $false, $false, $true, 1, $false, "57-29=28", 2) | Out-Null
$d.Content.Find.Execute("17+11=28", $true, $false, $false, $false, $false, $true, 1, $false, "26+44=70", 2) | Out-Null
$d.Content.Find.Execute("28+23=51", $true, $false, $false, $false, $false, $true, 1, $false, "54+21=75", 2) | Out-Null
$d.Content.Find.Execute("14+8=22", $true, $false, $false, $false, $false, $true, 1, $false, "51+0=51", 2) | Out-Null
$d.Content.Find.Execute("0+78=78", $true, $false, $false, $false, $false, $true, 1, $false, "96-27=69", 2) | Out-Null
$d.Content.Find.Execute("63+18=81", $true, $false, $false, $false, $false, $true, 1, $false, "46-34=12", 2) | Out-Null
$d.Content.Find.Execute("39+39=78", $true, $false, $false, $false, $false, $true, 1, $false, "89-71=18", 2) | Out-Null
$d.Content.Find.Execute("81-24=57", $true, $false, $false, $false, $false, $true, 1, $false, "69+11=80", 2) | Out-Null
$d.Content.Find.Execute("11+80=91", $true, $false, $false, $false, $false, $true, 1, $false, "8+28=36", 2) | Out-Null
$d.Content.Find.Execute("39-11=28", $true, $false, $false, $false, $false, $true, 1, $false, "48+50=98", 2) | Out-Null
$d.Content.Find.Execute("34+56=90", $true, $false, $false, $false, $false, $true, 1, $false, "86+11=97", 2) | Out-Null
$d.Content.Find.Execute("75-10=65", $true, $false, $false, $false, $false, $true, 1, $false, "52-41=11", 2) | Out-Null
$d.Content.Find.Execute("17+76=93", $true, $false, $false, $false, $false, $true, 1, $false, "49+25=74", 2) | Out-Null
$d.Content.Find.Execute("28+19=47", $true, $false, $false, $false, $false, $true, 1, $false, "40+10=50", 2) | Out-Null
$d.Content.Find.Execute("91-59=32", $true, $false, $false, $false, $false, $true, 1, $false, "55+39=94", 2) | Out-Null
$d.Content.Find.Execute("56-30=26", $true, $false, $false, $false, $false, $true, 1, $false, "33+13=46", 2) | Out-Null
$d.Content.Find.Execute("94-32=62", $true, $false, $false, $false, $false, $true, 1, $false, "16+68=84", 2) | Out-Null
$d.Content.Find.Execute("35-13=22", $true, $false, $false, $false, $false, $true, 1, $false, "94+3=97", 2) | Out-Null
$d.Content.Find.Execute("21-19=2", $true, $false, $false, $false, $false, $true, 1, $false, "49-23=26", 2) | Out-Null
$d.Content.Find.Execute("48+4=52", $true, $false, $false, $false, $false, $true, 1, $false, "17-4=13", 2) | Out-Null
$d.Content.Find.Execute("41+30=71", $true, $false, $false, $false, $false, $true, 1, $false, "89-12=77", 2) | Out-Null
$d.Content.Find.Execute("8+43=51", $true, $false, $false, $false, $false, $true, 1, $false, "15+13=28", 2) | Out-Null
$d.Content.Find.Execute("10-10=0", $true, $false, $false, $false, $false, $true, 1, $false, "42-33=9", 2) | Out-Null
$d.Content.Find.Execute("30+67=97", $true, $false, $false, $false, $false, $true, 1, $false, "45+22=67", 2) | Out-Null
$d.Content.Find.Execute("67-27=40", $true, $false, $false, $false, $false, $true, 1, $false, "30-21=9", 2) | Out-Null
$d.Content.Find.Execute("70-14=56", $true, $false, $false, $false, $false, $true, 1, $false, "41+0=41", 2) | Out-Null
$d.Content.Find.Execute("42+33=75", $true, $false, $false, $false, $false, $true, 1, $false, "78-19=59", 2) | Out-Null
$d.Content.Find.Execute("28+14=42", $true, $false, $false, $false, $false, $true, 1, $false, "58-6=52", 2) | Out-Null
$d.Content.Find.Execute("59-12=47", $true, $false, $false, $false, $false, $true, 1, $false, "82-68=14", 2) | Out-Null
$d.Content.Find.Execute("36+54=90", $true, $false, $false, $false, $false, $true, 1, $false, "59+38=97", 2) | Out-Null
$d.Content.Find.Execute("13+69=82", $true, $false, $false, $false, $false, $true, 1, $false, "67-50=17", 2) | Out-Null
$d.Content.Find.Execute("9+8=17", $true, $false, $false, $false, $false, $true, 1, $false, "80+16=96", 2) | Out-Null
$d.Content.Find.Execute("93-60=33", $true, $false, $false, $false, $false, $true, 1, $false, "85-18=67", 2) | Out-Null
$d.Content.Find.Execute("91-80=11", $true, $false, $false, $false, $false, $true, 1, $false, "95-6=89", 2) | Out-Null
$d.Content.Find.Execute("4+56=60", $true, $false, $false, $false, $false, $true, 1, $false, "39-7=32", 2) | Out-Null
$d.Content.Find.Execute("83-71=12", $true, $false, $false, $false, $false, $true, 1, $false, "78-19=59", 2) | Out-Null
$d.Content.Find.Execute("53+41=94", $true, $false, $false, $false, $false, $true, 1, $false, "17+12=29", 2) | Out-Null
$d.Content.Find.Execute("40+5=45", $true, $false, $false, $false, $false, $true, 1, $false, "42+32=74", 2) | Out-Null
$d.Content.Find.Execute("43-12=31", $true, $false, $false, $false, $false, $true, 1, $false, "65-42=23", 2) | Out-Null
$d.Content.Find.Execute("38+49=87", $true, $false, $false, $false, $false, $true, 1, $false, "2+79=81", 2) | Out-Null
$d.Content.Find.Execute("55+12=67", $true, $false, $false, $false, $false, $true, 1, $false, "6+66=72", 2) | Out-Null
$d.Content.Find.Execute("17+72=89", $true, $false, $false, $false, $false, $true, 1, $false, "0+81=81", 2) | Out-Null
$d.Content.Find.Execute("94-56=38", $true, $false, $false, $false, $false, $true, 1, $false, "35-15=20", 2) | Out-Null
$d.Content.Find.Execute("23+37=60", $true, $false, $false, $false, $false, $true, 1, $false, "59-4=55", 2) | Out-Null
$d.Content.Find.Execute("79-40=39", $true, $false, $false, $false, $false, $true, 1, $false, "36+18=54", 2) | Out-Null
$d.Content.Find.Execute("10+50=60", $true, $false, $false, $false, $false, $true, 1, $false, "86-27=59", 2) | Out-Null
$d.Content.Find.Execute("51-2=49", $true, $false, $false, $false, $false, $true, 1, $false, "61-45=16", 2) | Out-Null
$d.Content.Find.Execute("14+7=21", $true, $false, $false, $false, $false, $true, 1, $false, "47+18=65", 2) | Out-Null
$d.Content.Find.Execute("3+79=82", $true, $false, $false, $false, $false, $true, 1, $false, "49-46=3", 2) | Out-Null
$d.Content.Find.Execute("56-44=12", $true, $false, $false, $false, $false, $true, 1, $false, "4+47=51", 2) | Out-Null
$d.Content.Find.Execute("98-64=34", $true, $false, $false, $false, $false, $true, 1, $false, "54+41=95", 2) | Out-Null
$d.Content.Find.Execute("84+6=90", $true, $false, $false, $false, $false, $true, 1, $false, "57+36=93", 2) | Out-Null
$d.Content.Find.Execute("59-39=20", $true, $false, $false, $false, $false, $true, 1, $false, "19+25=44", 2) | Out-Null
$d.Content.Find.Execute("18-5=13", $true, $false, $false, $false, $false, $true, 1, $false, "49+47=96", 2) | Out-Null
$d.Content.Find.Execute("19+31=50", $true, $false, $false, $false, $false, $true, 1, $false, "98-77=21", 2) | Out-Null
$d.Content.Find.Execute("43+18=61", $true, $false, $false, $false, $false, $true, 1, $false, "49-18=31", 2) | Out-Null
$d.Content.Find.Execute("57-23=34", $true, $false, $false, $false, $false, $true, 1, $false, "12-1=11", 2) | Out-Null
$d.Content.Find.Execute("92-19=73", $true, $false, $false, $false, $false, $true, 1, $false, "52-33=19", 2) | Out-Null
$d.Content.Find.Execute("40+30=70", $true, $false, $false, $false, $false, $true, 1, $false, "8+12=20", 2) | Out-Null
$d.Content.Find.Execute("66+1=67", $true, $false, $false, $false, $false, $true, 1, $false, "72-47=25", 2) | Out-Null
$d.Content.Find.Execute("80-72=8", $true, $false, $false, $false, $false, $true, 1, $false, "55-30=25", 2) | Out-Null
$d.Content.Find.Execute("84-37=47", $true, $false, $false, $false, $false, $true, 1, $false, "5+87=92", 2) | Out-Null
$d.Content.Find.Execute("7+17=24", $true, $false, $false, $false, $false, $true, 1, $false, "58-23=35", 2) | Out-Null
$d.Content.Find.Execute("32-1=31", $true, $false, $false, $false, $false, $true, 1, $false, "57-26=31", 2) | Out-Null
$d.Content.Find.Execute("73+17=90", $true, $false, $false, $false, $false, $true, 1, $false, "81-27=54", 2) | Out-Null
$d.Content.Find.Execute("78-78=0", $true, $false, $false, $false, $false, $true, 1, $false, "65+17=82", 2) | Out-Null
$d.Content.Find.Execute("11+86=97", $true, $false, $false, $false, $false, $true, 1, $false, "43-9=34", 2) | Out-Null
$d.Content.Find.Execute("80-36=44", $true, $false, $false, $false, $false, $true, 1, $false, "84-13=71", 2) | Out-Null
$d.Content.Find.Execute("31+0=31", $true, $false, $false, $false, $false, $true, 1, $false, "26+63=89", 2) | Out-Null
$d.Content.Find.Execute("36+35=71", $true, $false, $false, $false, $false, $true, 1, $false, "92-15=77", 2) | Out-Null
